# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the K column (G2:G38) with recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0, 3, 9, 4, 6, 11, 6, 5, 8, 4, 3, 2, 3, 4, 9, 5, 10, 4, 1, 1, 2, 2, 5, 3, 5, 5, 11, 8, 6, 10, 5, 3, 5, 6, 2, 5, 3)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
